# ---------------------------------------------------------------------------
# Reproduces:
#   1. ppt/slides/slide6.xml : the table's <a:tableStyleId> changes from
#      {B29CD498-56E8-4941-B6F7-0F1205724FFF} to {9A412D6B-1B7C-488F-B518-B695C8F4913C}
#   2. ppt/theme/theme2.xml  : the deck's live theme colour scheme changes
#      from the "Integral" palette to the "Office Theme" palette (the slide
#      master / presentation theme). (ppt/theme/theme1.xml backs the Notes
#      Master's theme only, and is not reachable through the PowerPoint
#      object model - there is no VBA/COM surface that edits a Notes
#      Master's theme colours independently from the presentation theme.)
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on the table found on slide 6 --------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{9A412D6B-1B7C-488F-B518-B695C8F4913C}")
    }
}

# --- 2. Theme colour scheme: Integral -> Office Theme -----------------------
# Order of Colors(1..12): dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# RGB values are packed as COLORREF (BGR), matching VBA's RGB()/.RGB.
$firstSlide = $p.Slides.Item(1)
$colors = $firstSlide.ThemeColorScheme

$officeTheme = @{
    1  = 0           # dk1      000000
    2  = 16777215    # lt1      FFFFFF
    3  = 6968388      # dk2      44546A
    4  = 15132391     # lt2      E7E6E6
    5  = 13998939     # accent1  5B9BD5
    6  = 3243501       # accent2  ED7D31
    7  = 10855845      # accent3  A5A5A5
    8  = 49407         # accent4  FFC000
    9  = 12874308      # accent5  4472C4
    10 = 4697456        # accent6  70AD47
    11 = 12673797       # hlink    0563C1
    12 = 7491477        # folHlink 954F72
}

for ($i = 1; $i -le 12; $i++) {
    $colors.Colors($i).RGB = $officeTheme[$i]
}
